$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 2
$ws.Range("H2").Value = 1185.3334
$ws.Range("J2").Value = 3244.5
$ws.Range("L2").Value = 3244.5
$ws.Range("N2").Value = -3470.5

$ws = $wb.Worksheets.Item("ALC")  # row 62
$ws.Range("H62").Value = 4648.9287
$ws.Range("I62").Value = 4582.1665
$ws.Range("J62").Value = 4699
$ws.Range("K62").Value = 4582.1665
$ws.Range("L62").Value = 4699
$ws.Range("M62").Value = -3958.1665
$ws.Range("N62").Value = -5947

$ws = $wb.Worksheets.Item("ALC")  # row 65
$ws.Range("H65").Value = 4648.9287
$ws.Range("I65").Value = 4582.1665
$ws.Range("J65").Value = 4699
$ws.Range("K65").Value = 22910.8325
$ws.Range("L65").Value = 23495
$ws.Range("M65").Value = -19790.8325
$ws.Range("N65").Value = -29735

$ws = $wb.Worksheets.Item("ALC")  # row 92
$ws.Range("H92").Value = 408
$ws.Range("I92").Value = 464
$ws.Range("J92").Value = 228.8
$ws.Range("K92").Value = 464
$ws.Range("L92").Value = 228.8
$ws.Range("M92").Value = 784
$ws.Range("N92").Value = -2724.8

$ws = $wb.Worksheets.Item("ALC")  # row 138
$ws.Range("H138").Value = 4257.533
$ws.Range("I138").Value = 2758.6538
$ws.Range("K138").Value = 8275.9614
$ws.Range("M138").Value = -3135.9614

$ws = $wb.Worksheets.Item("ARM")  # row 32
$ws.Range("H32").Value = 4879.566
$ws.Range("I32").Value = 1899.5493
$ws.Range("K32").Value = 1899.5493
$ws.Range("M32").Value = -1612.5493

$ws = $wb.Worksheets.Item("ARM")  # row 61
$ws.Range("H61").Value = 8694.950999999999
$ws.Range("I61").Value = 7691.9414
$ws.Range("J61").Value = 13566.714
$ws.Range("K61").Value = 7691.9414
$ws.Range("L61").Value = 13566.714
$ws.Range("M61").Value = -7479.9414
$ws.Range("N61").Value = -13990.714

$ws = $wb.Worksheets.Item("ARM")  # row 122
$ws.Range("H122").Value = 4925.269
$ws.Range("I122").Value = 4946.1665
$ws.Range("J122").Value = 4674.5
$ws.Range("K122").Value = 14838.4995
$ws.Range("L122").Value = 14023.5
$ws.Range("M122").Value = -12388.4995
$ws.Range("N122").Value = -18923.5

$ws = $wb.Worksheets.Item("ARM")  # row 132
$ws.Range("H132").Value = 3305.257
$ws.Range("I132").Value = 1676.625
$ws.Range("J132").Value = 6858.636
$ws.Range("K132").Value = 5029.875
$ws.Range("L132").Value = 20575.908
$ws.Range("M132").Value = -2499.875
$ws.Range("N132").Value = -25635.908

$ws = $wb.Worksheets.Item("ARM")  # row 136
$ws.Range("H136").Value = 8694.950999999999
$ws.Range("I136").Value = 7691.9414
$ws.Range("J136").Value = 13566.714
$ws.Range("K136").Value = 23075.8242
$ws.Range("L136").Value = 40700.142
$ws.Range("M136").Value = -20525.8242
$ws.Range("N136").Value = -45800.142

$ws = $wb.Worksheets.Item("BSM")  # row 54
$ws.Range("H54").Value = 31293
$ws.Range("J54").Value = 40793.6
$ws.Range("L54").Value = 40793.6
$ws.Range("N54").Value = -41761.6

$ws = $wb.Worksheets.Item("BSM")  # row 58
$ws.Range("H58").Value = 19225
$ws.Range("J58").Value = 19225
$ws.Range("L58").Value = 19225
$ws.Range("N58").Value = -19813

$ws = $wb.Worksheets.Item("BSM")  # row 94
$ws.Range("H94").Value = 870.4583
$ws.Range("I94").Value = 571.0476
$ws.Range("K94").Value = 571.0476
$ws.Range("M94").Value = -120.0476

$ws = $wb.Worksheets.Item("BSM")  # row 134
$ws.Range("H134").Value = 4916.25
$ws.Range("I134").Value = 1999.4445
$ws.Range("K134").Value = 5998.333500000001
$ws.Range("M134").Value = -3463.333500000001

$ws = $wb.Worksheets.Item("CRP")  # row 58
$ws.Range("H58").Value = 6076.846
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 16333
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 16333
$ws.Range("M58").Value = -2797
$ws.Range("N58").Value = -16739

$ws = $wb.Worksheets.Item("CRP")  # row 99
$ws.Range("H99").Value = 2543.1428
$ws.Range("I99").Value = 1999.5
$ws.Range("J99").Value = 3902.25
$ws.Range("K99").Value = 1999.5
$ws.Range("L99").Value = 3902.25
$ws.Range("M99").Value = -501.5
$ws.Range("N99").Value = -6898.25

$ws = $wb.Worksheets.Item("CRP")  # row 126
$ws.Range("H126").Value = 2543.1428
$ws.Range("I126").Value = 1999.5
$ws.Range("J126").Value = 3902.25
$ws.Range("K126").Value = 5998.5
$ws.Range("L126").Value = 11706.75
$ws.Range("M126").Value = -3528.5
$ws.Range("N126").Value = -16646.75

$ws = $wb.Worksheets.Item("CRP")  # row 136
$ws.Range("H136").Value = 6076.846
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 16333
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 48999
$ws.Range("M136").Value = -6450
$ws.Range("N136").Value = -54099

$ws = $wb.Worksheets.Item("CUL")  # row 5
$ws.Range("H5").Value = 2401
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()  # remove cell

$ws = $wb.Worksheets.Item("CUL")  # row 12
$ws.Range("H12").Value = 192.3125
$ws.Range("J12").Value = 228.75
$ws.Range("L12").Value = 686.25
$ws.Range("N12").Value = -1032.25

$ws = $wb.Worksheets.Item("CUL")  # row 107
$ws.Range("H107").Value = 633.6
$ws.Range("J107").Value = 633.6
$ws.Range("L107").Value = 1900.8
$ws.Range("N107").Value = -5740.8

$ws = $wb.Worksheets.Item("CUL")  # row 113
$ws.Range("H113").Value = 1553.7778
$ws.Range("I113").Value = 503
$ws.Range("K113").Value = 1509
$ws.Range("M113").Value = 661  # new cell

$ws = $wb.Worksheets.Item("CUL")  # row 129
$ws.Range("H129").Value = 11920133
$ws.Range("I129").Value = 11201.2
$ws.Range("J129").Value = 41692460
$ws.Range("K129").Value = 33603.60000000001
$ws.Range("L129").Value = 125077380
$ws.Range("M129").Value = -28603.60000000001
$ws.Range("N129").Value = -125087380

$ws = $wb.Worksheets.Item("CUL")  # row 131
$ws.Range("H131").Value = 588129.5600000001
$ws.Range("I131").Value = 914.3333
$ws.Range("J131").Value = 1401196.9
$ws.Range("K131").Value = 2742.9999
$ws.Range("L131").Value = 4203590.699999999
$ws.Range("M131").Value = 2297.0001
$ws.Range("N131").Value = -4213670.699999999

$ws = $wb.Worksheets.Item("CUL")  # row 135
$ws.Range("H135").Value = 2401
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()  # remove cell

$ws = $wb.Worksheets.Item("GSM")  # row 43
$ws.Range("H43").Value = 42733.168
$ws.Range("I43").Value = 40000
$ws.Range("J43").Value = 43279.8
$ws.Range("K43").Value = 40000
$ws.Range("L43").Value = 43279.8
$ws.Range("M43").Value = -39849
$ws.Range("N43").Value = -43581.8

$ws = $wb.Worksheets.Item("GSM")  # row 122
$ws.Range("H122").Value = 5742.174
$ws.Range("I122").Value = 4389.421
$ws.Range("J122").Value = 12167.75
$ws.Range("K122").Value = 13168.263
$ws.Range("L122").Value = 36503.25
$ws.Range("M122").Value = -10718.263
$ws.Range("N122").Value = -41403.25

$ws = $wb.Worksheets.Item("LTW")  # row 7
$ws.Range("H7").Value = 50916.695
$ws.Range("I7").Value = 59797.633
$ws.Range("K7").Value = 59797.633
$ws.Range("M7").Value = -59685.633

$ws = $wb.Worksheets.Item("LTW")  # row 40
$ws.Range("H40").Value = 4518.923
$ws.Range("I40").Value = 2849
$ws.Range("K40").Value = 2849
$ws.Range("M40").Value = -2713

$ws = $wb.Worksheets.Item("LTW")  # row 122
$ws.Range("H122").Value = 7263.1816
$ws.Range("I122").Value = 7173.75
$ws.Range("J122").Value = 7501.6665
$ws.Range("K122").Value = 21521.25
$ws.Range("L122").Value = 22504.9995
$ws.Range("M122").Value = -19071.25
$ws.Range("N122").Value = -27404.9995

$ws = $wb.Worksheets.Item("LTW")  # row 126
$ws.Range("H126").Value = 50916.695
$ws.Range("I126").Value = 59797.633
$ws.Range("K126").Value = 179392.899
$ws.Range("M126").Value = -176922.899

$ws = $wb.Worksheets.Item("WVR")  # row 41
$ws.Range("H41").Value = 16887.75
$ws.Range("J41").Value = 16817
$ws.Range("L41").Value = 16817
$ws.Range("N41").Value = -17597

$ws = $wb.Worksheets.Item("WVR")  # row 81
$ws.Range("H81").Value = 1722.0435
$ws.Range("I81").Value = 1695.8096
$ws.Range("J81").Value = 1997.5
$ws.Range("K81").Value = 3391.6192
$ws.Range("L81").Value = 3995
$ws.Range("M81").Value = -2330.6192
$ws.Range("N81").Value = -6117

$ws = $wb.Worksheets.Item("WVR")  # row 84
$ws.Range("H84").Value = 1722.0435
$ws.Range("I84").Value = 1695.8096
$ws.Range("J84").Value = 1997.5
$ws.Range("K84").Value = 16958.096
$ws.Range("L84").Value = 19975
$ws.Range("M84").Value = -11654.096
$ws.Range("N84").Value = -30583

$ws = $wb.Worksheets.Item("WVR")  # row 122
$ws.Range("H122").Value = 3818.5293
$ws.Range("I122").Value = 2956.5386
$ws.Range("J122").Value = 6620
$ws.Range("K122").Value = 8869.6158
$ws.Range("L122").Value = 19860
$ws.Range("M122").Value = -6419.6158
$ws.Range("N122").Value = -24760

$ws = $wb.Worksheets.Item("WVR")  # row 126
$ws.Range("H126").Value = 1998.25
$ws.Range("I126").Value = 1855.1428
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 5565.428400000001
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -3095.428400000001
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("WVR")  # row 132
$ws.Range("H132").Value = 3715.5
$ws.Range("J132").Value = 20329.666
$ws.Range("L132").Value = 60988.99800000001
$ws.Range("N132").Value = -66048.99800000001

